$wb = $excel.ActiveWorkbook

# Column F ("想去人数" / "want to go" count) updates on sheet "展览" (sheet1)
$wsExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    5  = 15326
    6  = 414
    8  = 687
    9  = 15329
    11 = 8905
    14 = 76
    20 = 534
    27 = 68
    31 = 43
    32 = 35
    33 = 236
    34 = 293
    37 = 5441
}
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Same updates mirrored on sheet "全部类型" (sheet4), which repeats the same
# events at slightly different row offsets (extra rows inserted earlier).
$wsAll = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    5  = 15326
    6  = 414
    8  = 687
    9  = 15329
    11 = 8905
    15 = 76
    21 = 534
    28 = 68
    34 = 43
    35 = 35
    36 = 236
    37 = 293
    40 = 5441
}
foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
